# Update legacy GSC export data:
# - Drop the oldest date row (2025-10-22) from the "Chart" sheet, shifting
#   all subsequent rows up by one.
# - Append a new trailing row for the new date (2026-01-20) with zeroed
#   metric values, matching the existing trailing zero pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-10-22). Deleting it shifts every
# following row up by one, so row 2 becomes 2025-10-23, etc., and the
# data now ends at row 90 (2026-01-19).
$ws.Rows.Item(2).Delete()

# Append the new date row 91 for 2026-01-20 with the same zero values
# used for newly added days. The date column stores plain text (not a
# real date value), so write it with a leading apostrophe to keep Excel
# from auto-converting it to a date serial, then copy the plain-text
# formatting from an existing date cell so no extra style is introduced.
$ws.Range("A91").Value = "'2026-01-20"
$ws.Range("A90").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 0
$excel.CutCopyMode = 0
